$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.654.67"
$ws.Range("E2").Value = "  +4.04%  "

$ws.Range("D3").Value = "1.699.15"
$ws.Range("E3").Value = "  +2.47%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'318.32"
$ws.Range("E5").Value = "  +3.35%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "'0.3965"
$ws.Range("E7").Value = "  +2.16%  "

$ws.Range("D8").Value = "'0.4038"
$ws.Range("E8").Value = "  +2.32%  "

$ws.Range("E9").Value = "  +9.98%  "

$ws.Range("D10").Value = "'54.29"
$ws.Range("E10").Value = "  +9.92%  "

$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("D12").Value = "'0.08823"
$ws.Range("E12").Value = "  +2.04%  "

$ws.Range("D13").Value = "'7.284"
$ws.Range("E13").Value = "  +8.14%  "

$ws.Range("D14").Value = "'23.40"
$ws.Range("E14").Value = "  +3.19%  "

$ws.Range("E15").Value = "  +1.94%  "

$ws.Range("D16").Value = "'7.651"
$ws.Range("E16").Value = "  +6.19%  "

$ws.Range("D17").Value = "1.701.01"
$ws.Range("E17").Value = "  +2.29%  "

$ws.Range("D18").Value = "'101.39"
$ws.Range("E18").Value = "  +1.72%  "

$ws.Range("D19").Value = "'0.07093"
$ws.Range("E19").Value = "  +4.83%  "

$ws.Range("D20").Value = "'19.79"
$ws.Range("E20").Value = "  +4.30%  "

$ws.Range("D21").Value = "'6.902"
$ws.Range("E21").Value = "  +3.91%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").Value = "'14.15"
$ws.Range("E23").Value = "  +2.80%  "

$ws.Range("D24").Value = "24.641.84"
$ws.Range("E24").Value = "  +4.03%  "

$ws.Range("D25").Value = "'3.072"
$ws.Range("E25").Value = "  +11.51%  "

$ws.Range("D26").Value = "'2.331"
$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").Value = "'22.48"
$ws.Range("E27").Value = "  +3.64%  "

$ws.Range("D28").Value = "'159.95"
$ws.Range("E28").Value = "  +1.89%  "

$ws.Range("D29").Value = "'5.249"
$ws.Range("E29").Value = "  +1.84%  "

$ws.Range("D30").Value = "'134.68"
$ws.Range("E30").Value = "  +4.22%  "

$ws.Range("D31").Value = "'7.505"
$ws.Range("E31").Value = "  +16.00%  "

$ws.Range("D32").Value = "'1.120"
$ws.Range("E32").Value = "  -0.69%  "

$ws.Range("D33").Value = "1.890.70"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("D34").Value = "'7.500"
$ws.Range("E34").Value = "  +16.75%  "

$ws.Range("D35").Value = "'0.08597"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").Value = "'11.56"
$ws.Range("E36").Value = "  +12.62%  "

$ws.Range("E37").Value = "  +4.53%  "

$ws.Range("D38").Value = "'1.954"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("D39").Value = "'14.75"
$ws.Range("E39").Value = "  +3.05%  "

$ws.Range("D40").Value = "'0.02806"
$ws.Range("E40").Value = "  +11.45%  "

$ws.Range("D41").Value = "'0.09074"
$ws.Range("E41").Value = "  +3.69%  "

$ws.Range("D42").Value = "'0.7782"
$ws.Range("E42").Value = "  +3.76%  "

$ws.Range("D43").Value = "'1.468"
$ws.Range("E43").Value = "  +1.22%  "

$ws.Range("D44").Value = "'0.7269"
$ws.Range("E44").Value = "  +3.95%  "

$ws.Range("D45").Value = "'15.64"
$ws.Range("E45").Value = "  +4.75%  "

$ws.Range("D46").Value = "'2.532"
$ws.Range("E46").Value = "  +6.75%  "

$ws.Range("D47").Value = "'4.233"
$ws.Range("E47").Value = "  +4.10%  "

$ws.Range("D48").Value = "'1.369"
$ws.Range("E48").Value = "  +16.05%  "

$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").Value = "'141.71"
$ws.Range("E50").Value = "  +2.77%  "

$ws.Range("D51").Value = "'0.08064"
$ws.Range("E51").Value = "  +4.25%  "
